# Move Emon conccn & sconcss from pre identified missing to pre ignored file #384.
# These two variable rows (Emon/conccn and Emon/sconcss) together with the two
# blank rows immediately following them are removed from the worksheet, which
# shifts every subsequent row up by four positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 715 (Emon conccn) and 716 (Emon sconcss) hold the two variable
# definitions being relocated out of this workbook; rows 717-718 are the
# blank rows that sat between this block and the next ("Omon talkos") entry.
# Deleting the whole 715:718 block removes the data and shifts the remaining
# rows up by four, matching the new row numbering.
$ws.Range("A715:A718").EntireRow.Delete()

# Restore the view/selection state to what it was after the edit.
$ws.Range("A714").Select()
$excel.ActiveWindow.ScrollRow = 699
$excel.ActiveWindow.DisplayGridlines = $true
